$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Insert a new empty paragraph right after the heading paragraph, then
    # give it the ListBullet style and the new docente's text.
    $target.Range.InsertParagraphAfter()
    $newPara = $target.Next()
    $newPara.Style = "ListBullet"
    $newPara.Range.Text = "6712818 - Mauricio Lamano Ferreira"
}
